$d = $word.ActiveDocument

# --- Locate the paragraph containing "5) I am still trying to deal with ..." ---
$searchRange = $d.Content
$searchRange.Find.Execute("5) I am still trying to deal with", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$fivePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $searchRange.Start -and $p.Range.End -ge $searchRange.End) {
        $fivePara = $p
        break
    }
}

if ($fivePara -eq $null) {
    throw "Could not locate the '5) I am still trying...' paragraph"
}

# The paragraph immediately before it is the blank spacer paragraph that
# should be removed along with it.
$fiveIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $fivePara.Range.Start) {
        $fiveIndex = $i
        break
    }
}
$spacerPara = $d.Paragraphs.Item($fiveIndex - 1)

# --- Remove the "5) ..." paragraph (its text and its paragraph mark). ---
$fivePara.Range.Delete() | Out-Null

# --- Remove the now-adjacent blank spacer paragraph (and its mark too). ---
$spacerPara.Range.Delete() | Out-Null

# --- Find the _GoBack bookmark paragraph that now directly follows, and
#     insert a brand-new empty paragraph right after it. Use a half-open
#     [Start, End) containment test so the lookup is unambiguous even
#     when the bookmark sits exactly on a paragraph boundary. ---
$bmStart = $d.Bookmarks("_GoBack").Range.Start
$bookmarkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $bmStart -and $bmStart -lt $p.Range.End) {
        $bookmarkPara = $p
        break
    }
}

if ($bookmarkPara -eq $null) {
    throw "Could not locate the _GoBack bookmark paragraph"
}

$insPoint = $d.Range($bookmarkPara.Range.End, $bookmarkPara.Range.End)
$insPoint.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>") | Out-Null
